$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.222383260726929
$ws.Range("B1").Value = 2.296833753585815
$ws.Range("C1").Value = 3.343996524810791
$ws.Range("D1").Value = 2.125917673110962
$ws.Range("E1").Value = 1.313796162605286
